$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new/changed cell text values in the same order the shared strings
# table picks up new entries (A13, D16, A12, A17, A18, A16, A19, A20).

# Row 13: A13 -> new "Từ tác vụ đầu tiên của Xuân Sang" follow-up text
$ws.Range("A13").Value = "- Từ tác vụ đầu tiên của Xuân Sang`n    + Viết thêm hàm đọc/ghi một người thuê từ file/ra file`n[Bổ sung]`n+ Viết thêm hàm thayDoiThongTinTaiKhoanMatKhau()"

# Row 16: D16 -> "Xuân Lam"
$ws.Range("D16").Value = "Xuân Lam"

# Row 12: A12 -> ThueMotMay struct text, now prefixed with a new bracket line
$ws.Range("A12").Value = "[Tạo mới và viết vào file ThueMotMay.hpp]`n- Viết code khai báo biến ThueMotMay với struct, thông tin mỗi máy tính như sau:`n    + Tài khoản (char[30])`n    + Mật khẩu (char[30])`n    + Số điện thoại (char[15])`n    + Con trỏ maytinh trỏ đến máy`n- Kiểm tra số máy có hợp lệ không`n- Chọn máy`n- Nhập một người thuê`n- Xuất một người thuê theo chiều ngang/dọc`n`n[Bổ sung]`n+ Đã viết thêm hàm layViTriCuaMay() vào file NhieuMayTinh.hpp để gọi vào hàm chonMay()`n"

# Row 17: A17 -> new "Từ tác vụ đầu tiên của Xuân Lam" task, D17 -> Hải Sơn
$ws.Range("A17").Value = "- Từ tác vụ đầu tiên của Xuân Lam, viết thêm hàm`n+ Kiểm tra tài khoản đã tồn tại hay chưa`n+ Thêm một người thuê vào danh sách`n+ Xuất danh sách người thuê"
$ws.Range("D17").Value = "Hải Sơn"

# Row 18: A18 -> ThueNhieuMay.hpp edit task, D18 -> Quốc Thắng
$ws.Range("A18").Value = "[Viết thêm hàm vào file ThueNhieuMay.hpp]`n- Chỉnh sửa danh sách người thuê, kiểm tra số điện thoại có tồn tại hay không nếu có thì yêu cầu nhập thông tin tài khoản và mật khẩu mới.`n- Xóa một người thuê khỏi danh sách`n- Đọc danh sách người thuê`n- Ghi danh sách người thuê`n[Bổ sung]`n+ Viết thêm hàm kiểm tra tài khoản, số điện thoại đã tồn tại hay chưa"
$ws.Range("D18").Value = "Quốc Thắng"

# Row 16: A16 -> ThueNhieuMay.hpp linked-list task
$ws.Range("A16").Value = "[Tạo mới và viết vào file ThueNhieuMay.hpp]`n- Khai báo dữ liệu danh sách liên kết đơn với  dữ liệu lưu trữ là 1 khách thuê (1 phần tử)`n- Viết các hàm:`n+ createList();`n+ createNode();`n+ addNodeInTail();`n+ removeNodeInHead();`n+ removeNodeInTail();`n+ giaiPhongVungNhoDanhSachNguoiThue();"

# Row 19: A19 -> ThueNhieuMay.hpp write-to-file task, D19 -> Xuân Sang
$ws.Range("A19").Value = "[Viết thêm hàm vào file ThueNhieuMay.hpp]`n- Ghi danh sách người thuê ra file`n- Xuất danh sách người thuê ra file`n- Tìm kiếm người thuê theo tài khoản`n"
$ws.Range("D19").Value = "Xuân Sang"

# Row 20: A20 -> switch-case editing task, D20 -> Thanh Sang
$ws.Range("A20").Value = "[Viết thêm hàm vào file ThueNhieuMay.hpp]`n- Tạo một switch case để nhập thông tin muốn chỉnh sửa. VD: 1. Thay đổi tài khoản 2. Thay đổi mật khẩu 3. Thay đổi số điện thoại…"
$ws.Range("D20").Value = "Thanh Sang"

# Row heights (matches real-Excel autofit/manual resize after the edits)
$ws.Rows(12).RowHeight = 288.75
$ws.Rows(13).RowHeight = 137.45
$ws.Rows(14).RowHeight = 49.15
$ws.Rows(15).RowHeight = 70.15
$ws.Rows(16).RowHeight = 201
$ws.Rows(17).RowHeight = 83.25
$ws.Rows(18).RowHeight = 200.25
$ws.Rows(19).RowHeight = 81.75
$ws.Rows(20).RowHeight = 79.5

# Sheet view updates: zoom + selection
$excel.ActiveWindow.Zoom = 60
$ws.Range("D15").Select()

